$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1574.1818
$ws.Range("I98").Value = 1884.1333
$ws.Range("J98").Value = 910
$ws.Range("K98").Value = 1884.1333
$ws.Range("L98").Value = 910
$ws.Range("M98").Value = -386.1333
$ws.Range("N98").Value = -3906

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2734.5
$ws.Range("I106").Value = 2734.5
$ws.Range("K106").Value = 2734.5
$ws.Range("M106").Value = -2103.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1888.5652
$ws.Range("I107").Value = 1493.0588
$ws.Range("J107").Value = 3009.1667
$ws.Range("K107").Value = 1493.0588
$ws.Range("L107").Value = 3009.1667
$ws.Range("M107").Value = 426.9412
$ws.Range("N107").Value = -6849.1667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 20003054
$ws.Range("J113").Value = 3603.75
$ws.Range("L113").Value = 3603.75
$ws.Range("N113").Value = -10111.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1574.1818
$ws.Range("I122").Value = 1884.1333
$ws.Range("J122").Value = 910
$ws.Range("K122").Value = 5652.3999
$ws.Range("L122").Value = 2730
$ws.Range("M122").Value = -3202.3999
$ws.Range("N122").Value = -7630

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 38000
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 38000
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 38000
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -47800

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 144.75
$ws.Range("I135").Value = 133.7
$ws.Range("K135").Value = 1203.3
$ws.Range("M135").Value = 1331.7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1597.5416
$ws.Range("I137").Value = 1292.05
$ws.Range("J137").Value = 3125
$ws.Range("K137").Value = 3876.15
$ws.Range("L137").Value = 9375
$ws.Range("M137").Value = -1326.15
$ws.Range("N137").Value = -14475

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2722.4768
$ws.Range("I32").Value = 2830.1865
$ws.Range("K32").Value = 2830.1865
$ws.Range("M32").Value = -2543.1865

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 45000
$ws.Range("I62").Value = 45000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 45000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -44376
$ws.Range("N62").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 45000
$ws.Range("I65").Value = 45000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 135000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -131880
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 1687666.6
$ws.Range("J92").Value = 1687666.6
$ws.Range("L92").Value = 1687666.6
$ws.Range("N92").Value = -1692658.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1210.6086
$ws.Range("I110").Value = 1037.7727
$ws.Range("K110").Value = 1037.7727
$ws.Range("M110").Value = 1007.2273

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2304.653
$ws.Range("I132").Value = 2022.55
$ws.Range("K132").Value = 6067.65
$ws.Range("M132").Value = -3537.65

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8333974.5
$ws.Range("I94").Value = 11364134
$ws.Range("K94").Value = 11364134
$ws.Range("M94").Value = -11363683

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1929.4
$ws.Range("I107").Value = 1683.2858
$ws.Range("J107").Value = 2503.6667
$ws.Range("K107").Value = 1683.2858
$ws.Range("L107").Value = 2503.6667
$ws.Range("M107").Value = 236.7141999999999
$ws.Range("N107").Value = -6343.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 877.6727
$ws.Range("I31").Value = 722.03125
$ws.Range("J31").Value = 1094.2174
$ws.Range("K31").Value = 722.03125
$ws.Range("L31").Value = 1094.2174
$ws.Range("M31").Value = -427.03125
$ws.Range("N31").Value = -1684.2174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 877.6727
$ws.Range("I34").Value = 722.03125
$ws.Range("J34").Value = 1094.2174
$ws.Range("K34").Value = 722.03125
$ws.Range("L34").Value = 1094.2174
$ws.Range("M34").Value = -520.03125
$ws.Range("N34").Value = -1498.2174

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 28000
$ws.Range("J50").Value = 28000
$ws.Range("L50").Value = 28000
$ws.Range("N50").Value = -29250

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 19110.889
$ws.Range("J88").Value = 20874.75
$ws.Range("L88").Value = 20874.75
$ws.Range("N88").Value = -21686.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 19110.889
$ws.Range("J91").Value = 20874.75
$ws.Range("L91").Value = 20874.75
$ws.Range("N91").Value = -23682.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9084.6
$ws.Range("I132").Value = 10105.917
$ws.Range("K132").Value = 30317.751
$ws.Range("M132").Value = -27787.751

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 9525006
$ws.Range("I134").Value = 11495352
$ws.Range("K134").Value = 34486056
$ws.Range("M134").Value = -34483521

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2000
$ws.Range("J75").Value = 2000
$ws.Range("L75").Value = 6000
$ws.Range("N75").Value = -7996

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H78").Value = 2000
$ws.Range("J78").Value = 2000
$ws.Range("L78").Value = 18000
$ws.Range("N78").Value = -27984

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 7957.7393
$ws.Range("J96").Value = 7957.7393
$ws.Range("L96").Value = 23873.2179
$ws.Range("N96").Value = -27991.2179

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 6339.6665
$ws.Range("I107").Value = 581.625
$ws.Range("J107").Value = 10946.1
$ws.Range("K107").Value = 1744.875
$ws.Range("L107").Value = 32838.3
$ws.Range("M107").Value = 175.125
$ws.Range("N107").Value = -36678.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 637.55554
$ws.Range("I113").Value = 514.2857
$ws.Range("J113").Value = 680.7
$ws.Range("K113").Value = 1542.8571
$ws.Range("L113").Value = 2042.1
$ws.Range("M113").Value = 627.1428999999998
$ws.Range("N113").Value = -6382.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 643.0714
$ws.Range("I122").Value = 475
$ws.Range("J122").Value = 867.1667
$ws.Range("K122").Value = 4275
$ws.Range("L122").Value = 7804.5003
$ws.Range("M122").Value = -1825
$ws.Range("N122").Value = -12704.5003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 13890059
$ws.Range("I129").Value = 41667496
$ws.Range("J129").Value = 3789172.8
$ws.Range("K129").Value = 125002488
$ws.Range("L129").Value = 11367518.4
$ws.Range("M129").Value = -124997488
$ws.Range("N129").Value = -11377518.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17545244
$ws.Range("J131").Value = 1502.9387
$ws.Range("L131").Value = 4508.8161
$ws.Range("N131").Value = -14588.8161

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 19569718
$ws.Range("I70").Value = 20836956
$ws.Range("J70").Value = 18187276
$ws.Range("K70").Value = 20836956
$ws.Range("L70").Value = 18187276
$ws.Range("M70").Value = -20836686
$ws.Range("N70").Value = -18187816

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 19569718
$ws.Range("I73").Value = 20836956
$ws.Range("J73").Value = 18187276
$ws.Range("K73").Value = 20836956
$ws.Range("L73").Value = 18187276
$ws.Range("M73").Value = -20836020
$ws.Range("N73").Value = -18189148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3493.9412
$ws.Range("I126").Value = 2075
$ws.Range("J126").Value = 3930.5386
$ws.Range("K126").Value = 6225
$ws.Range("L126").Value = 11791.6158
$ws.Range("M126").Value = -3755
$ws.Range("N126").Value = -16731.6158

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2195
$ws.Range("I132").Value = 1715.7333
$ws.Range("J132").Value = 4591.3335
$ws.Range("K132").Value = 5147.199900000001
$ws.Range("L132").Value = 13774.0005
$ws.Range("M132").Value = -2617.199900000001
$ws.Range("N132").Value = -18834.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2346.75
$ws.Range("I61").Value = 2344
$ws.Range("J61").Value = 2349.5
$ws.Range("K61").Value = 2344
$ws.Range("L61").Value = 2349.5
$ws.Range("M61").Value = -2142
$ws.Range("N61").Value = -2753.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 15000
$ws.Range("J64").Value = 15000
$ws.Range("L64").Value = 15000
$ws.Range("N64").Value = -15450

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H67").Value = 15000
$ws.Range("J67").Value = 15000
$ws.Range("L67").Value = 15000
$ws.Range("N67").Value = -16560

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2346.75
$ws.Range("I113").Value = 2344
$ws.Range("J113").Value = 2349.5
$ws.Range("K113").Value = 2344
$ws.Range("L113").Value = 2349.5
$ws.Range("M113").Value = -174
$ws.Range("N113").Value = -6689.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 56678668
$ws.Range("I122").Value = 56678668
$ws.Range("K122").Value = 170036004
$ws.Range("M122").Value = -170033554

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6752
$ws.Range("I136").Value = 11263.4
$ws.Range("J136").Value = 1739.3334
$ws.Range("K136").Value = 33790.2
$ws.Range("L136").Value = 5218.0002
$ws.Range("M136").Value = -31240.2
$ws.Range("N136").Value = -10318.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 49502.8
$ws.Range("J123").Value = 49502.8
$ws.Range("L123").Value = 49502.8
$ws.Range("N123").Value = -59302.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2548.342
$ws.Range("I132").Value = 2166.1936
$ws.Range("K132").Value = 6498.5808
$ws.Range("M132").Value = -3968.5808

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 575.4545000000001
$ws.Range("I136").Value = 388.125
$ws.Range("J136").Value = 1075
$ws.Range("K136").Value = 1164.375
$ws.Range("L136").Value = 3225
$ws.Range("M136").Value = 1385.625
$ws.Range("N136").Value = -8325

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 44715
$ws.Range("J141").Value = 44715
$ws.Range("L141").Value = 44715
$ws.Range("N141").Value = -55075
